$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New hangman words appended to the list (in the order they were added)
$newWords = @(
    "xanadu",
    "cowabunga",
    "wannabe",
    "radical",
    "cool",
    "duh",
    "cd",
    "nerd",
    "dude",
    "chill",
    "whatever",
    "dynamite",
    "gnarly ",
    "disco",
    "groove"
)

# Find the first empty row below the existing data in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row = $lastRow + 1

foreach ($word in $newWords) {
    $ws.Cells.Item($row, 1).Value = $word
    $row = $row + 1
}

$lastRow = $row - 1

# Sort the word list (A2:A<lastRow>) alphabetically, ascending
$sortRange = $ws.Range("A2:A$lastRow")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortRange, 0, 1, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Restore the view/selection state left after the edits
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B20").Select()
